$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift every timestamp in column A (rows 2-97) forward from 2024-06-25 to
#    2024-07-30 (the time-of-day portion is unchanged, only the date moves).
# ---------------------------------------------------------------------------
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    $updated = $current -replace "^2024-06-25", "2024-07-30"
    $cell.Value = $updated
}

# ---------------------------------------------------------------------------
# 2) Re-split the weather-derived figures across locations for both the base
#    and intraday models (columns B..F, rows 24-49).
# ---------------------------------------------------------------------------
$newData = @{
    24 = @(0, 0, 0, 0, 0)
    25 = @(0, 0, 0, 0, 0)
    26 = @(1642.31844329834, 0.0016423184432983, 0, 0.0008211592216491, 0.0002052898054122)
    27 = @(12806.68037923177, 0.0128066803792317, 0.0016423184432983, 0.007224499411265, 0.0018061248528162)
    28 = @(27043.73956298828, 0.0270437395629882, 0.0128066803792317, 0.0199252099711099, 0.0049813024927774)
    29 = @(54038.58296712239, 0.0540385829671223, 0.0270437395629882, 0.0405411612650552, 0.0101352903162637)
    30 = @(88996.98274739584, 0.0889969827473958, 0.0540385829671223, 0.0715177828572591, 0.0178794457143147)
    31 = @(147492.6302083333, 0.1474926302083333, 0.0889969827473958, 0.1182448064778645, 0.0295612016194661)
    32 = @(271976.7798665364, 0.2719767798665364, 0.1474926302083333, 0.2097347050374348, 0.0524336762593586)
    33 = @(429977.9078776041, 0.4299779078776041, 0.2719767798665364, 0.3509773438720702, 0.08774433596801751)
    34 = @(633277.2262369791, 0.6332772262369791, 0.4299779078776041, 0.5316275670572916, 0.1329068917643228)
    35 = @(863057.2936197917, 0.8630572936197917, 0.6332772262369791, 0.7481672599283853, 0.1870418149820962)
    36 = @(1109280.420572917, 1.109280420572917, 0.8630572936197917, 0.9861688570963542, 0.2465422142740885)
    37 = @(1343670.201171875, 1.343670201171875, 1.109280420572917, 1.226475310872396, 0.3066188277180988)
    38 = @(1593599.343098958, 1.593599343098958, 1.343670201171875, 1.468634772135417, 0.3671586930338541)
    39 = @(1845423.793619792, 1.845423793619792, 1.593599343098958, 1.719511568359375, 0.4298778920898437)
    40 = @(2070057.217447917, 2.070057217447917, 1.845423793619792, 1.957740505533854, 0.4894351263834635)
    41 = @(2285382.825520833, 2.285382825520833, 2.070057217447917, 2.177720021484375, 0.5444300053710937)
    42 = @(2497296.0078125, 2.4972960078125, 2.285382825520833, 2.391339416666667, 0.5978348541666667)
    43 = @(2689228.684895833, 2.689228684895833, 2.4972960078125, 2.593262346354167, 0.6483155865885417)
    44 = @(2875311.48828125, 2.87531148828125, 2.689228684895833, 2.782270086588541, 0.6955675216471353)
    45 = @(2976534.640625, 2.976534640625, 2.87531148828125, 2.925923064453125, 0.7314807661132812)
    46 = @(0, 0, 2.976534640625, 1.4882673203125, 0.372066830078125)
    47 = @(0, 0, 0, 0, 0)
    48 = @(0, 0, 0, 0, 0)
    49 = @(0, 0, 0, 0, 0)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}
